$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (first sheet) ---
$ws1 = $wb.Worksheets.Item(1)

# Version: 5.0.0 -> 6.0.0
$ws1.Cells.Item(3, 2).Value = "6.0.0"

# Date: update publish date
$ws1.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty -> now "Alvearie Team"
$ws1.Cells.Item(9, 2).Value = "Alvearie Team"

# Remove the duplicated "Contact" row (row 11), leaving just one row at 10
$ws1.Rows.Item(11).Delete()

# Row 10 becomes "Jurisdiction" / "United States of America"
$ws1.Cells.Item(10, 1).Value = "Jurisdiction"
$ws1.Cells.Item(10, 2).Value = "United States of America"

# --- Sheet "Elements" (second sheet) ---
$ws2 = $wb.Worksheets.Item(2)

# Row 2 (root Extension element): Short / Definition now reflect the attribution code title/description
$ws2.Cells.Item(2, 11).Value = "Attribution Code"
$ws2.Cells.Item(2, 12).Value = "FFV initiative attribution method code"
